$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("searchProduct")

# Replace old "Woo Album #1" in A8 and add the new album rows first so the
# shared-string table gains the "Woo Album N" entries before "Woo Logo".
$ws.Range("A8").Value = "Woo Album 1"
$ws.Range("A9").Value = "Woo Album 2"
$ws.Range("A10").Value = "Woo Album 3"
$ws.Range("A11").Value = "Woo Album 4"

# Shift remaining original content down: A5 keeps "Patient Ninja", A6 keeps
# "Premium Quality", A7 keeps "Ship Your Idea" (unchanged values, already there)
$ws.Range("A5").Value = "Patient Ninja"
$ws.Range("A6").Value = "Premium Quality"
$ws.Range("A7").Value = "Ship Your Idea"

# Update existing cell A4 (was "Ninja Silhouette" -> now "Woo Logo") last so
# this new shared string is appended at the end of the table.
$ws.Range("A4").Value = "Woo Logo"

$ws.Range("B7").Select()
